$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 2023 (row 19) Q1 figure and recompute the yearly total (column F).
$ws.Range("B19").Value = 343676.77580000012
$ws.Range("F19").Value = 1498823.5917000002

# Add the 2024* (row 20) Q2 figure and recompute its running total (column F).
$ws.Range("C20").Value = 512712.72640000097
$ws.Range("F20").Value = 638444.07430000114

# Leave the cursor on L10, matching the saved selection in the workbook.
$ws.Range("L10").Select()
